$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/28/2025  Through  8/3/2025"

# --- Weekly crime-stat table updates (rows 15-31) ---
$ws.Cells.Item(15,7).NumberFormat = "@"
$ws.Cells.Item(15,7).Value = "0"
$ws.Cells.Item(15,8).NumberFormat = "@"
$ws.Cells.Item(15,8).Value = "***.*"
$ws.Cells.Item(15,12).Value = -9.090909090909
$ws.Cells.Item(15,13).Value = 25
$ws.Cells.Item(16,3).Value = 1
$ws.Cells.Item(16,5).Value = -50
$ws.Cells.Item(16,6).Value = 12
$ws.Cells.Item(16,7).Value = 9
$ws.Cells.Item(16,8).Value = 33.333333333333
$ws.Cells.Item(16,9).Value = 64
$ws.Cells.Item(16,10).Value = 71
$ws.Cells.Item(16,11).Value = -9.859154929577
$ws.Cells.Item(16,12).Value = -17.948717948717
$ws.Cells.Item(16,13).Value = -57.615894039735
$ws.Cells.Item(16,14).Value = -88.191881918819
$ws.Cells.Item(17,3).Value = 7
$ws.Cells.Item(17,4).Value = 4
$ws.Cells.Item(17,5).Value = 75
$ws.Cells.Item(17,6).Value = 12
$ws.Cells.Item(17,7).Value = 20
$ws.Cells.Item(17,8).Value = -40
$ws.Cells.Item(17,9).Value = 130
$ws.Cells.Item(17,10).Value = 169
$ws.Cells.Item(17,11).Value = -23.076923076923
$ws.Cells.Item(17,12).Value = 0
$ws.Cells.Item(17,13).Value = 120.338983050847
$ws.Cells.Item(17,14).Value = -7.801418439716
$ws.Cells.Item(18,3).Value = 1
$ws.Cells.Item(18,4).Value = 3
$ws.Cells.Item(18,5).Value = -66.666666666666
$ws.Cells.Item(18,6).Value = 13
$ws.Cells.Item(18,7).Value = 18
$ws.Cells.Item(18,8).Value = -27.777777777777
$ws.Cells.Item(18,9).Value = 124
$ws.Cells.Item(18,10).Value = 129
$ws.Cells.Item(18,11).Value = -3.875968992248
$ws.Cells.Item(18,12).Value = -25.301204819277
$ws.Cells.Item(18,13).Value = -20
$ws.Cells.Item(18,14).Value = -86.036036036036
$ws.Cells.Item(19,3).Value = 11
$ws.Cells.Item(19,4).Value = 15
$ws.Cells.Item(19,5).Value = -26.666666666666
$ws.Cells.Item(19,6).Value = 37
$ws.Cells.Item(19,7).Value = 39
$ws.Cells.Item(19,8).Value = -5.128205128205
$ws.Cells.Item(19,9).Value = 315
$ws.Cells.Item(19,10).Value = 336
$ws.Cells.Item(19,11).Value = -6.25
$ws.Cells.Item(19,12).Value = -9.221902017291
$ws.Cells.Item(19,13).Value = 17.100371747211
$ws.Cells.Item(19,14).Value = -12.983425414364
$ws.Cells.Item(20,3).Value = 3
$ws.Cells.Item(20,4).Value = 7
$ws.Cells.Item(20,5).Value = -57.142857142857
$ws.Cells.Item(20,6).Value = 23
$ws.Cells.Item(20,7).Value = 29
$ws.Cells.Item(20,8).Value = -20.689655172413
$ws.Cells.Item(20,9).Value = 200
$ws.Cells.Item(20,10).Value = 209
$ws.Cells.Item(20,11).Value = -4.306220095693
$ws.Cells.Item(20,12).Value = -4.761904761904
$ws.Cells.Item(20,13).Value = 66.666666666666
$ws.Cells.Item(20,14).Value = -93.045897079276
$ws.Cells.Item(21,3).Value = 23
$ws.Cells.Item(21,4).Value = 31
$ws.Cells.Item(21,5).Value = -25.806451612903
$ws.Cells.Item(21,6).Value = 97
$ws.Cells.Item(21,7).Value = 115
$ws.Cells.Item(21,8).Value = -15.652173913043
$ws.Cells.Item(21,9).Value = 846
$ws.Cells.Item(21,10).Value = 929
$ws.Cells.Item(21,11).Value = -8.934337997847
$ws.Cells.Item(21,12).Value = -10.191082802547
$ws.Cells.Item(21,13).Value = 10.588235294117
$ws.Cells.Item(21,14).Value = -82.498965659909
$ws.Cells.Item(22,3).NumberFormat = "#,##0"
$ws.Cells.Item(22,3).Value = 1
$ws.Cells.Item(22,4).Value = 1
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 1
$ws.Cells.Item(22,8).Value = -75
$ws.Cells.Item(22,9).Value = 15
$ws.Cells.Item(22,10).Value = 24
$ws.Cells.Item(22,11).Value = -37.5
$ws.Cells.Item(22,12).Value = -40
$ws.Cells.Item(22,13).Value = 0
$ws.Cells.Item(23,3).Value = 2
$ws.Cells.Item(23,4).NumberFormat = "#,##0"
$ws.Cells.Item(23,4).Value = 1
$ws.Cells.Item(23,5).NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Cells.Item(23,5).Value = 100
$ws.Cells.Item(23,6).Value = 5
$ws.Cells.Item(23,7).Value = 2
$ws.Cells.Item(23,8).Value = 150
$ws.Cells.Item(23,9).Value = 34
$ws.Cells.Item(23,10).Value = 33
$ws.Cells.Item(23,11).Value = 3.030303030303
$ws.Cells.Item(23,12).Value = -29.166666666666
$ws.Cells.Item(23,13).Value = 61.904761904761
$ws.Cells.Item(24,3).Value = 35
$ws.Cells.Item(24,4).Value = 24
$ws.Cells.Item(24,5).Value = 45.833333333333
$ws.Cells.Item(24,6).Value = 88
$ws.Cells.Item(24,7).Value = 83
$ws.Cells.Item(24,8).Value = 6.024096385542
$ws.Cells.Item(24,9).Value = 644
$ws.Cells.Item(24,10).Value = 699
$ws.Cells.Item(24,11).Value = -7.868383404864
$ws.Cells.Item(24,12).Value = -23.515439429928
$ws.Cells.Item(24,13).Value = 11.805555555555
$ws.Cells.Item(25,3).Value = 8
$ws.Cells.Item(25,4).Value = 16
$ws.Cells.Item(25,5).Value = -50
$ws.Cells.Item(25,7).Value = 31
$ws.Cells.Item(25,8).Value = -29.032258064516
$ws.Cells.Item(25,9).Value = 232
$ws.Cells.Item(25,10).Value = 301
$ws.Cells.Item(25,11).Value = -22.923588039867
$ws.Cells.Item(25,12).Value = -22.666666666666
$ws.Cells.Item(26,3).Value = 7
$ws.Cells.Item(26,4).Value = 16
$ws.Cells.Item(26,5).Value = -56.25
$ws.Cells.Item(26,6).Value = 42
$ws.Cells.Item(26,7).Value = 47
$ws.Cells.Item(26,8).Value = -10.63829787234
$ws.Cells.Item(26,9).Value = 317
$ws.Cells.Item(26,10).Value = 304
$ws.Cells.Item(26,11).Value = 4.276315789473
$ws.Cells.Item(26,12).Value = 12.41134751773
$ws.Cells.Item(26,13).Value = 23.346303501945
$ws.Cells.Item(27,10).Value = 19
$ws.Cells.Item(27,11).Value = -36.842105263157
$ws.Cells.Item(27,12).Value = -33.333333333333
$ws.Cells.Item(28,3).NumberFormat = "@"
$ws.Cells.Item(28,3).Value = "0"
$ws.Cells.Item(28,4).NumberFormat = "#,##0"
$ws.Cells.Item(28,4).Value = 1
$ws.Cells.Item(28,5).NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Cells.Item(28,5).Value = -100
$ws.Cells.Item(28,6).Value = 3
$ws.Cells.Item(28,7).Value = 2
$ws.Cells.Item(28,8).Value = 50
$ws.Cells.Item(28,10).Value = 24
$ws.Cells.Item(28,11).Value = 58.333333333333
$ws.Cells.Item(28,12).Value = 46.153846153846
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "0"
$ws.Cells.Item(31,5).NumberFormat = "@"
$ws.Cells.Item(31,5).Value = "***.*"
$ws.Cells.Item(31,7).Value = 2
